$d = $word.ActiveDocument

$d.Content.Find.Execute("94×22=2068", $true, $false, $false, $false, $false, $true, 1, $false, "90×89=8010", 2)
$d.Content.Find.Execute("56×19=1064", $true, $false, $false, $false, $false, $true, 1, $false, "81×72=5832", 2)
$d.Content.Find.Execute("25×72=1800", $true, $false, $false, $false, $false, $true, 1, $false, "73×51=3723", 2)
$d.Content.Find.Execute("93×62=5766", $true, $false, $false, $false, $false, $true, 1, $false, "14×69=966", 2)
$d.Content.Find.Execute("15×99=1485", $true, $false, $false, $false, $false, $true, 1, $false, "64×80=5120", 2)
$d.Content.Find.Execute("28×47=1316", $true, $false, $false, $false, $false, $true, 1, $false, "77×22=1694", 2)
$d.Content.Find.Execute("65×66=4290", $true, $false, $false, $false, $false, $true, 1, $false, "33×73=2409", 2)
$d.Content.Find.Execute("92×15=1380", $true, $false, $false, $false, $false, $true, 1, $false, "18×59=1062", 2)
$d.Content.Find.Execute("46×67=3082", $true, $false, $false, $false, $false, $true, 1, $false, "29×82=2378", 2)
$d.Content.Find.Execute("93×81=7533", $true, $false, $false, $false, $false, $true, 1, $false, "69×76=5244", 2)
$d.Content.Find.Execute("77×19=1463", $true, $false, $false, $false, $false, $true, 1, $false, "97×80=7760", 2)
$d.Content.Find.Execute("31×76=2356", $true, $false, $false, $false, $false, $true, 1, $false, "49×42=2058", 2)
$d.Content.Find.Execute("88×86=7568", $true, $false, $false, $false, $false, $true, 1, $false, "93×12=1116", 2)
$d.Content.Find.Execute("26×44=1144", $true, $false, $false, $false, $false, $true, 1, $false, "35×52=1820", 2)
$d.Content.Find.Execute("90×45=4050", $true, $false, $false, $false, $false, $true, 1, $false, "18×51=918", 2)
$d.Content.Find.Execute("61×31=1891", $true, $false, $false, $false, $false, $true, 1, $false, "67×16=1072", 2)
$d.Content.Find.Execute("83×63=5229", $true, $false, $false, $false, $false, $true, 1, $false, "95×78=7410", 2)
$d.Content.Find.Execute("89×91=8099", $true, $false, $false, $false, $false, $true, 1, $false, "16×84=1344", 2)
$d.Content.Find.Execute("57×45=2565", $true, $false, $false, $false, $false, $true, 1, $false, "92×14=1288", 2)
$d.Content.Find.Execute("67×22=1474", $true, $false, $false, $false, $false, $true, 1, $false, "22×63=1386", 2)
$d.Content.Find.Execute("93×88=8184", $true, $false, $false, $false, $false, $true, 1, $false, "61×73=4453", 2)
$d.Content.Find.Execute("24×31=744", $true, $false, $false, $false, $false, $true, 1, $false, "67×85=5695", 2)
$d.Content.Find.Execute("33×29=957", $true, $false, $false, $false, $false, $true, 1, $false, "74×36=2664", 2)
$d.Content.Find.Execute("26×38=988", $true, $false, $false, $false, $false, $true, 1, $false, "64×36=2304", 2)
$d.Content.Find.Execute("67×77=5159", $true, $false, $false, $false, $false, $true, 1, $false, "80×77=6160", 2)

Write-Host "All replacements applied"
